# Update "想去人数" (interest count) figures in F column on both the
# "展览" and "全部类型" sheets. Same row->value map applies to each sheet.
$wb = $excel.ActiveWorkbook

$updates = @{
    8  = 2056
    11 = 4479
    16 = 127
    20 = 3332
    25 = 83
    29 = 61
    33 = 1977
    34 = 356
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
